$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Pow10($mantissa, $exp) {
    return $mantissa * [Math]::Pow(10, $exp)
}

# Row 2 updates
$ws.Range("E2").Value = 24.66000000000042
$ws.Range("H2").Value = Pow10 1.364329369739056 -16
$ws.Range("K2").Value = 57.59772675983201
$ws.Range("L2").Value = "[49.29350404214529, 65.90194947751874]"
$ws.Range("O2").Value = 1.402552876377425
$ws.Range("P2").Value = "[1.2641844311742707, 1.540921321580579]"
$ws.Range("S2").Value = 58.75131015945689
$ws.Range("T2").Value = "[53.750092460192924, 63.752527858720846]"
$ws.Range("W2").Value = 19.15531531531564
$ws.Range("X2").Value = 18.61225225225257
$ws.Range("Y2").Value = 19.69837837837871

# Row 3 updates
$ws.Range("E3").Value = 25.49000000000055
$ws.Range("H3").Value = Pow10 1.364329369739056 -16
$ws.Range("K3").Value = 60.02794618592772
$ws.Range("L3").Value = "[51.49677053962682, 68.55912183222861]"
$ws.Range("O3").Value = 2.03150035457358
$ws.Range("P3").Value = "[1.8931319093704255, 2.169868799776734]"
$ws.Range("S3").Value = 58.51806041973698
$ws.Range("T3").Value = "[53.97281136688414, 63.06330947258983]"
$ws.Range("W3").Value = 17.24848848848886
$ws.Range("X3").Value = 16.68714714714751
$ws.Range("Y3").Value = 17.80982982983021
